# Fix position names ("Coordinator" -> "Coordinators") on the Organizers sheet,
# and update the matching bio text that repeated the same typo.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Organizers")

# Column D holds each organizer's "position" string. Every row whose position
# ends in "Coordinator" (singular) should be pluralized to "Coordinators".
for ($r = 2; $r -le 24; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $current = $cell.Value()
    if ($current -ne $null -and $current -like "*Coordinator") {
        $cell.Value = $current + "s"
    }
}

# Vishvak's bio (G18) repeats the same "Coordinator" position name inline;
# fix it the same way so it stays consistent with the corrected position.
$bioCell = $ws.Range("G18")
$bio = $bioCell.Value()
$bioCell.Value = $bio -replace "HackUTD Industry Coordinator ", "HackUTD Industry Coordinators "

# Move the saved cell cursor from D15 to D16, matching the workbook's new
# selection state.
$ws.Range("D16").Select()
